# Diagramas de robustez y secuencia CU 10 y 11
# Mark the related tasks as "Hecho" (done) and log 1 hour of work on
# "Día 7" for the CU 10 / CU 11 rows (25 y 26), along with the related
# minor touch-ups on CU 22 / CU 08 / CU 06 / CU 09 rows (23, 24, 27, 28).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Rows 23-24 (CU 22 / CU 08) were already in progress: just flip the
# status to "Hecho".
$ws.Range("F23").Value = "Hecho"
$ws.Range("F24").Value = "Hecho"

# Rows 25-28 (CU 10, CU 11, CU 06, CU 09): mark done, record the
# estimated hour and the hour consumed on "Día 7" (column Z).
foreach ($r in 25..28) {
    $ws.Cells.Item($r, 6).Value = "Hecho"      # F: Estatus
    $ws.Cells.Item($r, 7).Value = 1            # G: Horas estimadas totales
    $ws.Cells.Item($r, 26).Value = 1           # Z: Día 7 - Cons.
}

# Move the active selection to where the user was working (Día 7 column,
# row 28).
$ws.Range("Z28").Select()
